# V2 of Overcuts — replace the captured-data table with the new dataset
# (12 rows, distances 0..1100) and retitle the two "Tile" labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 0. Stash copies of the two cell formats we need to reuse later, in a
#    scratch cell far outside the table, since the originals (E3 / E2)
#    get overwritten below. We clear the scratch cell again at the end.
# ---------------------------------------------------------------------
$ws.Range("E3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)     # Z1 = "highlighted" (red/bold) style
$ws.Range("E2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)     # Z2 = plain numeric style

$highlightFmt = $ws.Range("Z1")
$plainFmt     = $ws.Range("Z2")

# ---------------------------------------------------------------------
# 1. New table contents (row 2 .. row 13)
# ---------------------------------------------------------------------
$data = @(
    @{ Row=2;  A=0;    B=1222.55; C=1219.55; D=3;    E=3;    F="12in"; Hi=$false }
    @{ Row=3;  A=100;  B=1223.36; C=1219.79; D=3.57; E=3.57; F="12in"; Hi=$false }
    @{ Row=4;  A=200;  B=1224.34; C=1219.89; D=4.46; E=4.45; F="O8in"; Hi=$true  }
    @{ Row=5;  A=300;  B=1224.84; C=1219.99; D=4.86; E=4.85; F="O8in"; Hi=$true  }
    @{ Row=6;  A=400;  B=1225.53; C=1220.09; D=5.44; E=5.44; F="O8in"; Hi=$false }
    @{ Row=7;  A=500;  B=1225.09; C=1220.19; D=4.9;  E=4.9;  F="O8in"; Hi=$false }
    @{ Row=8;  A=600;  B=1225.55; C=1220.29; D=5.26; E=5.26; F="O8in"; Hi=$false }
    @{ Row=9;  A=700;  B=1225.82; C=1220.39; D=5.43; E=5.43; F="O8in"; Hi=$false }
    @{ Row=10; A=800;  B=1225.49; C=1220.49; D=5;    E=5;    F="O8in"; Hi=$false }
    @{ Row=11; A=900;  B=1225.47; C=1220.59; D=4.89; E=4.88; F="O8in"; Hi=$true  }
    @{ Row=12; A=1000; B=1225.5;  C=1220.69; D=4.81; E=4.81; F="O8in"; Hi=$false }
    @{ Row=13; A=1100; B=1225.58; C=1220.79; D=4.79; E=4.79; F="O8in"; Hi=$false }
)

foreach ($d in $data) {
    $r = $d.Row

    # Make sure A:E carry the plain numeric format (rows 7-13 are brand
    # new cells with no format yet; rows 2-6 already have it, this is a
    # harmless no-op refresh for those).
    $plainFmt.Copy()
    $ws.Range("A$r`:E$r").PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $d.A
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E

    if ($d.Hi) {
        $highlightFmt.Copy()
        $ws.Range("E$r").PasteSpecial(-4122)
        $ws.Cells.Item($r, 5).Value = $d.E
    }

    # Tile label text (keeps trailing CRLF like the source workbook).
    $ws.Range("F$r").Characters().Text = "$($d.F)`r`n"
}

# ---------------------------------------------------------------------
# 2. Clean up scratch cells
# ---------------------------------------------------------------------
$ws.Range("Z1").Clear()
$ws.Range("Z2").Clear()
